# A new weekly price-report row for "Espinaca" (Región Metropolitana) was
# inserted into the daily consolidated sheet. Excel shifts every existing
# row from 675 downward by one (old row 675 -> new row 676, ..., old row
# 761 -> new row 762) and the freshly inserted row 675 carries the new
# week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 675; this pushes rows 675-761 down to 676-762
# and bumps the sheet's used range from R761 to R762, exactly like Excel's
# own "Insert Sheet Rows" command.
$ws.Rows.Item(675).Insert()

# Populate the newly inserted row 675 with the new observation.
$ws.Range("A675").Value2 = 6
$ws.Range("B675").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C675").Value2 = 'Metropolitana'
$ws.Range("D675").Value2 = 45077
$ws.Range("E675").Value2 = 13
$ws.Range("F675").Value2 = 100112012
$ws.Range("G675").Value2 = 'Espinaca'
$ws.Range("H675").Value2 = 'Sin especificar'
$ws.Range("I675").Value2 = 'Primera'
$ws.Range("J675").Value2 = 580
$ws.Range("K675").Value2 = 4500
$ws.Range("L675").Value2 = 5000
$ws.Range("M675").Value2 = 4716
$ws.Range("N675").Value2 = '$/cuna 10 kilos'
$ws.Range("O675").Value2 = 'Región Metropolitana'
$ws.Range("P675").Value2 = 472
$ws.Range("Q675").Value2 = 10
$ws.Range("R675").Value2 = 'Hortaliza'
